# Switch the deck's design theme from the pink/violet "Integral" theme to
# the default blue/orange "Office Theme" palette (Design tab -> Office Theme).
#
# The slide master's theme (ppt/theme/theme1.xml) carries the 12 DrawingML
# theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Re-point
# each one at the stock "Office" theme color values via the slide's
# ThemeColorScheme, which edits the master theme in place (font scheme /
# format scheme are already the shared "Office" ones, so only the color
# scheme needs to change).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
